# Update cryptos list (Price and Volume(1h) columns) per GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.622.06"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "1.879.67"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").Value = "'316.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "'0.3934"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").Value = "'0.08414"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").Value = "'41.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "'6.280"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.37%  "
$ws.Range("D13").Value = "1.879.17"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "'7.282"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "'1.011"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'0.00001108"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "'91.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'0.06742"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'5.973"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "28.631.67"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "'2.249"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "2.096.74"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'162.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("D28").Value = "'20.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "'2.373"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "'127.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'1.059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("D33").Value = "'5.820"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'3.625"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'0.02462"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "'0.06552"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").Value = "'8.915"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "'1.268"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.39%  "
$ws.Range("D40").Value = "'1.198"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("D41").Value = "'0.6475"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").Value = "'11.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "'1.008"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'0.6078"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'13.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").Value = "'3.707"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("D48").Value = "'2.034"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'1.220"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").Value = "'122.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").Value = "'1.196"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.57%  "
